$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header of first column
$ws.Range("A1").Value = "MaxFES"

# Update the "MaxFES" column values (rows 2-14)
$maxfesValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $maxfesValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $maxfesValues[$i]
}

# Remove the "Run 50" column (column BA) entirely; this shifts the
# "Mean" column (BA) left into AZ.
$ws.Range("BA:BA").Delete() | Out-Null

# The deleted column leaves AZ1 still labeled "Run 50" (Excel does not
# auto-relabel header text); rename it to "Mean" to reflect the shifted
# column, and update the recomputed values underneath it.
$ws.Cells.Item(1, 52).Value = "Mean"
$azValues = @(140.61248675, 136.87602705, 107.72306394, 58.88029654, 50.43574943, 45.15122097, 41.63215018, 37.71878081, 33.71326577, 30.37423798, 26.93256095, 23.23273233, 20.02241539)
for ($i = 0; $i -lt $azValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $azValues[$i]
}
